$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new row at position 4 (pushes "Gardens on Spring Creek" and the
# rows below it down by one, keeping the alphabetical ordering of
# organizations: Colorado Water Center, Environmental Learning Center,
# [new] Fort Collins Museum of Discovery, Gardens on Spring Creek, ...)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the Fort Collins Museum of
# Discovery entry.
$ws.Range("A4").Value = "Fort Collins Museum of Discovery"
$ws.Range("B4").Value = "Museum"
$ws.Range("C4").Value = "Science museum"
$ws.Range("D4").Value = "Hands on displays and programs"
$ws.Range("E4").Value = "https://fcmod.org/"
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = -105.07795
$ws.Range("H4").Value = 40.5936

# Rebuild the Website-column hyperlinks since the row insert does not shift
# the existing hyperlink anchors automatically. Note row 6 (Houston Gardens)
# intentionally has no hyperlink, matching the source data (its cell is
# styled like a hyperlink but was never wired up to one).
$ws.Range("E2:E13").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), "https://watercenter.colostate.edu/")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://warnercnr.colostate.edu/elc/")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://fcmod.org/")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.fcgov.com/gardens/")
$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.northernwater.org/AboutUs/WaterEducation.aspx")
$ws.Hyperlinks.Add($ws.Range("E8"), "http://openwaterfoundation.org/")
$ws.Hyperlinks.Add($ws.Range("E9"), "https://www.poudreheritage.org/")
$ws.Hyperlinks.Add($ws.Range("E10"), "https://poudrelearningcenter.org/")
$ws.Hyperlinks.Add($ws.Range("E11"), "https://southplattebasin.com/")
$ws.Hyperlinks.Add($ws.Range("E12"), "https://www.watereducationcolorado.org/")
$ws.Hyperlinks.Add($ws.Range("E13"), "https://www.wgcd.org/conservation-education/")

# Restore the hyperlink cell style (Hyperlinks.Add changes formatting) on
# every Website cell, including E6 which keeps the style without a link.
$ws.Range("E2:E13").Style = "Hyperlink"

# Reset the view: no frozen/scrolled left column, and the active selection
# moved from E17 to G18.
$ws.Activate()
$ws.Range("G18").Select()

Write-Output "done"
